$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate header row and data rows (also updates dimension / shared strings) ---
$ws.Range("A1").Value = "Tags"
$ws.Range("B1").Value = "Description"
$ws.Range("A2").Value = "[FRAME2]"
$ws.Range("A3").Value = "[FRAME3]"
$ws.Range("A4").Value = "[FRAME1]"
$ws.Range("A5").Value = "[FAIL]"
$ws.Range("A6").Value = "[MAC]"
$ws.Range("A7").Value = "[CK]"

# --- Header formatting: bold font, centered alignment ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108

$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = -4108

# --- Column widths (Excel stores widths snapped to whole-pixel increments for the
# sheet's font, so these are the closest settable values that round-trip to the
# target stored widths of 22.28515625 and 21 characters) ---
$ws.Columns("A").ColumnWidth = 21.5
$ws.Columns("B").ColumnWidth = 20.15

# --- Selection state ---
$ws.Range("A9").Select()
